$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

function Replace-WholeWord($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND (whole word): $old"
    }
}

# Finds $old once, then appends each of the given strings as separate runs
# immediately following it (each InsertAfter + Collapse(0) step creates a
# new run boundary so the resulting OOXML keeps sentence/period runs split,
# matching how Word itself authored the original document).
function Append-Runs($old, [string[]]$pieces) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND (append-runs): $old"
        return
    }
    foreach ($piece in $pieces) {
        $r.InsertAfter($piece)
        $r.Collapse(0)
    }
}

# --- Title / author / email ---
Replace-Text "The Enduring Legacy of Marie Curie" "A Journey into the Realm of Mathematics: The Beauty of Patterns and Numbers"
Replace-Text "Anna Kowalski" "Mary Githinji"
Replace-WholeWord "anna" "mary"
Replace-Text "kowalski@emailworld" "githinji@schoolnet"
Replace-WholeWord "com" "edu"

# --- Paragraph 1 (intro) ---
Replace-Text "In the annals of science, few names resonate with such enduring respect and admiration as that of Marie Curie" "The world around us is permeated with numbers and patterns, from the gentle ripples of waves to the intricate design of snowflakes"

Replace-Text " A true pioneer in the fields of radioactivity and nuclear physics, her groundbreaking discoveries not only revolutionized our understanding of the atom but also paved the way for countless advancements in medicine, technology, and our overall comprehension of the universe" " Mathematics, often perceived as an abstract and enigmatic discipline, is the key to unlocking the underlying order and structure of these patterns"

Replace-Text " Curie's life and work stand as a testament to the transformative power of scientific inquiry, the indomitable spirit of human curiosity, and the boundless possibilities that lie at the intersection of perseverance and brilliance" " Delving into the realm of mathematics offers a journey of wonder and discovery, inviting us to explore the captivating beauty and inherent elegance of numbers"

Replace-Text "From her humble beginnings in Warsaw, Poland, Curie's journey into the realm of science was fueled by an insatiable thirst for knowledge and an unwavering determination to unravel the mysteries of nature" "In the mosaic of our universe, numbers hold a profound significance, shaping our perception of time, space, and quantity"

Replace-Text " Despite the formidable obstacles she faced as a woman in the male-dominated scientific world, Curie pressed on tirelessly, eventually securing her place as one of the most influential scientists of all time" " They are the threads that weave together the tapestry of reality, providing a framework for understanding the intricate dance of the cosmos"

Append-Runs " They are the threads that weave together the tapestry of reality, providing a framework for understanding the intricate dance of the cosmos" @(
    ".",
    " Mathematics enables us to decode this dance, revealing the underlying principles that govern the interactions of objects and phenomena",
    ".",
    " It grants us the power to unravel the mysteries of nature and create technologies that redefine the boundaries of human possibility"
)

# --- Paragraph 1 continued ---
Replace-Text "Throughout her illustrious career, Curie's unwavering dedication to her work led to a series of remarkable achievements" "Mathematics is not merely a set of abstract symbols or formulas; it is a living, dynamic discipline that permeates every aspect of our existence"

Replace-Text " Her pioneering research on radioactivity, conducted alongside her husband, Pierre Curie, resulted in the discovery of two new elements: polonium and radium" " It is the language of engineering marvels, guiding the construction of towering skyscrapers and sleek aircraft"

Replace-Text " These elements, with their extraordinary properties, would later find invaluable applications in medicine, including the development of groundbreaking cancer treatments" " It is the engine driving scientific discoveries, propelling humanity toward an ever-deepening understanding of the universe"

Append-Runs " It is the engine driving scientific discoveries, propelling humanity toward an ever-deepening understanding of the universe" @(
    ".",
    " And, equally importantly, mathematics is a source of inspiration and beauty, enticing us to marvel at the intricate patterns and symmetries that surround us"
)

# --- Summary paragraph ---
Replace-Text "Marie Curie's contributions to science are immeasurable" "The journey into the realm of mathematics is an invitation to discover the beauty and elegance of numbers and patterns"

Replace-Text " Her groundbreaking discoveries in radioactivity, including the identification of polonium and radium, revolutionized our understanding of the atom and laid the foundation for modern nuclear physics" " Mathematics is the key to unlocking the underlying order and structure of the universe, providing a framework for understanding the intricate dance of the cosmos"

Replace-Text " Curie's pioneering spirit, unwavering dedication, and tireless pursuit of knowledge serve as an inspiration to generations of scientists and continue to shape the course of scientific advancements to this day" " It has profound implications for our perception of time, space, and quantity, shaping our understanding of reality"

Append-Runs " It has profound implications for our perception of time, space, and quantity, shaping our understanding of reality" @(
    ".",
    " Moreover, mathematics is a vital tool for scientific discovery and ",
    "technological advancement, driving innovation and reshaping the world around us",
    ".",
    " It is a discipline that not only inspires and fascinates but also empowers us to unravel the mysteries of nature and create technologies that redefine the boundaries of human possibility"
)

# --- Trailing empty paragraph at end of document ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Host "DONE"
